$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the existing "Late" column (column N),
# shifting Late / heading(Date) / Outstanding one column to the right
# (N->O, O->P, P->Q). Matches column width of the preceding column (M).
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with Q13 selected -
# this also clears tabSelected on whichever sheet was previously active.
$ws.Activate()
$ws.Range("Q13").Select() | Out-Null
